$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7172284644194756
$ws1.Range("C2").Value = 0.8493975903614458
$ws1.Range("D2").Value = 0.5280898876404494
$ws1.Range("E2").Value = 0.651270207852194
$ws1.Range("F2").Value = 0.5713128038897893
$ws1.Range("G2").Value = 0.5358865662914779
$ws1.Range("H2").Value = 0.7172284644194756
$ws1.Range("I2").Value = 282
$ws1.Range("J2").Value = 50
$ws1.Range("K2").Value = 484
$ws1.Range("L2").Value = 252

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.6576086956521739
$ws2.Range("C2").Value = 0.9063670411985019
$ws2.Range("D2").Value = 0.7622047244094489

$ws2.Range("B3").Value = 0.8493975903614458
$ws2.Range("C3").Value = 0.5280898876404494
$ws2.Range("D3").Value = 0.651270207852194

$ws2.Range("B4").Value = 0.7172284644194756
$ws2.Range("C4").Value = 0.7172284644194756
$ws2.Range("D4").Value = 0.7172284644194756
$ws2.Range("E4").Value = 0.7172284644194756

$ws2.Range("B5").Value = 0.7535031430068099
$ws2.Range("C5").Value = 0.7172284644194756
$ws2.Range("D5").Value = 0.7067374661308214

$ws2.Range("B6").Value = 0.7535031430068099
$ws2.Range("C6").Value = 0.7172284644194756
$ws2.Range("D6").Value = 0.7067374661308214

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 484
$ws3.Range("C2").Value = 50
$ws3.Range("B3").Value = 252
$ws3.Range("C3").Value = 282
